$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price records for rows 3 and 4 were swapped (row 3 now holds
# what used to be row 4's data, and vice versa) for columns D, J, K, L, M, P.

$ws.Range("D3").Value = 44547
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 13500
$ws.Range("P3").Value = 750

$ws.Range("D4").Value = 44568
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("P4").Value = 861
